# Update wordpress_site_info.xlsx with latest site information

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update Scan Date ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "2025-08-05 21:29:02"

# --- Plugins sheet: update plugin version numbers (column B) ---
# Values that Excel would otherwise auto-convert to a number or a date
# (e.g. "1.14", "9.3.10") are prefixed with a leading apostrophe so they
# are stored as text, matching the original report's Version column.
$plugins = $wb.Worksheets.Item("Plugins")

$plugins.Cells.Item(8, 2).Value  = "'1.14"      # Admin Menu Editor
$plugins.Cells.Item(9, 2).Value  = "6.4.3"      # Advanced Custom Fields
$plugins.Cells.Item(10, 2).Value = "6.4.3"      # Advanced Custom Fields PRO
$plugins.Cells.Item(13, 2).Value = "'7.97"      # All-in-One WP Migration
$plugins.Cells.Item(14, 2).Value = "'2.73"      # All-in-One WP Migration Unlimited Extension
$plugins.Cells.Item(16, 2).Value = "'9.3.10"    # Solid Security Basic
$plugins.Cells.Item(17, 2).Value = "2.4.6"      # Broken Link Checker
$plugins.Cells.Item(18, 2).Value = "'6.1"       # Contact Form 7
$plugins.Cells.Item(20, 2).Value = "1.18.0"     # Custom Post Type UI
$plugins.Cells.Item(21, 2).Value = "1.21.1"     # GTM4WP
$plugins.Cells.Item(25, 2).Value = "9.6.1"      # Google Analytics for WordPress by MonsterInsights
$plugins.Cells.Item(26, 2).Value = "'1.4.15"    # Prisna GWT - Google Website Translator
$plugins.Cells.Item(27, 2).Value = "'2.9.14"    # Gravity Forms
$plugins.Cells.Item(30, 2).Value = "2.2.9"      # WPCode Lite
$plugins.Cells.Item(32, 2).Value = "4.9.2"      # PDF Embedder
$plugins.Cells.Item(34, 2).Value = "4.0.0"      # Post Type Switcher
$plugins.Cells.Item(35, 2).Value = "9.4.2"      # Really Simple Security
$plugins.Cells.Item(36, 2).Value = "9.4.1"      # Really Simple Security Pro
$plugins.Cells.Item(37, 2).Value = "5.5.2"      # Redirection
$plugins.Cells.Item(42, 2).Value = "3.1.3"      # TablePress
$plugins.Cells.Item(46, 2).Value = "2.25.6.26"  # UpdraftPlus - Backup/Restore
$plugins.Cells.Item(47, 2).Value = "4.64.5"     # User Role Editor Pro
$plugins.Cells.Item(48, 2).Value = "'25.6"      # Yoast SEO
$plugins.Cells.Item(49, 2).Value = "4.11.4"     # WP All Import Pro
$plugins.Cells.Item(50, 2).Value = "1.3.8"      # WP Fastest Cache
$plugins.Cells.Item(52, 2).Value = "2.7.3"      # WP Migrate Lite
$plugins.Cells.Item(53, 2).Value = "3.20.0"     # Smush
